$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 3.5
$ws.Range("J5").Value = 2.75
$ws.Range("K5").Value = 2.05

$ws.Range("S5").Value = 2.1
$ws.Range("T5").Value = 1.7
$ws.Range("W5").Value = 3.75
$ws.Range("X5").Value = 1.25

$ws.Range("AI5").Value = 9
$ws.Range("AM5").Value = 9.5
$ws.Range("AN5").Value = 17
$ws.Range("AS5").Value = 351
